# Update cryptocurrency price/volume data per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.991.20'
$ws.Range('E2').Value = '  +7.21%  '
$ws.Range('D3').Value = '1.744.21'
$ws.Range('E3').Value = '  +5.57%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = "'228.17"
$ws.Range('E5').Value = '  +4.28%  '
$ws.Range('D6').Value = "'0.5435"
$ws.Range('E6').Value = '  +3.72%  '
$ws.Range('D7').Value = "'1.002"
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').Value = "'0.2765"
$ws.Range('E8').Value = '  +4.00%  '
$ws.Range('D9').Value = "'0.06743"
$ws.Range('E9').Value = '  +6.34%  '
$ws.Range('D10').Value = "'21.69"
$ws.Range('E10').Value = '  +5.17%  '
$ws.Range('D11').Value = "'0.07781"
$ws.Range('E11').Value = '  +1.09%  '
$ws.Range('D12').Value = "'4.707"
$ws.Range('E12').Value = '  +2.07%  '
$ws.Range('D13').Value = '1.746.47'
$ws.Range('E13').Value = '  +4.88%  '
$ws.Range('D14').Value = '1.982.11'
$ws.Range('E14').Value = '  +5.45%  '
$ws.Range('D15').Value = "'0.5980"
$ws.Range('E15').Value = '  +6.63%  '
$ws.Range('D16').Value = '0.0₅8384'
$ws.Range('E16').Value = '  +2.22%  '
$ws.Range('D17').Value = "'68.94"
$ws.Range('E17').Value = '  +5.65%  '
$ws.Range('D18').Value = '27.977.03'
$ws.Range('D19').Value = "'224.44"
$ws.Range('E19').Value = '  +17.51%  '
$ws.Range('D20').Value = "'4.843"
$ws.Range('E20').Value = '  +3.04%  '
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').Value = "'10.91"
$ws.Range('E22').Value = '  +5.17%  '
$ws.Range('D23').Value = "'6.239"
$ws.Range('E23').Value = '  +4.26%  '
$ws.Range('D25').Value = "'146.27"
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('D26').Value = "'0.1246"
$ws.Range('E26').Value = '  +3.53%  '
$ws.Range('D27').Value = "'17.32"
$ws.Range('E27').Value = '  +8.85%  '
$ws.Range('D28').Value = "'1.671"
$ws.Range('E28').Value = '  +11.07%  '
$ws.Range('D29').Value = "'7.455"
$ws.Range('E29').Value = '  +2.75%  '
$ws.Range('D30').Value = "'0.05643"
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('D31').Value = "'1.314"
$ws.Range('E31').Value = '  +3.09%  '
$ws.Range('D32').Value = "'3.707"
$ws.Range('E32').Value = '  +6.11%  '
$ws.Range('D33').Value = "'3.520"
$ws.Range('E33').Value = '  +4.38%  '
$ws.Range('D34').Value = "'1.681"
$ws.Range('E34').Value = '  +6.38%  '
$ws.Range('D35').Value = "'0.9820"
$ws.Range('E35').Value = '  +3.34%  '
$ws.Range('D36').Value = "'2.859"
$ws.Range('E36').Value = '  +2.25%  '
$ws.Range('D37').Value = "'2.450"
$ws.Range('E37').Value = '  +1.64%  '
$ws.Range('D38').Value = "'0.5969"
$ws.Range('E38').Value = '  +3.71%  '
$ws.Range('D39').Value = "'0.01665"
$ws.Range('E39').Value = '  +4.46%  '
$ws.Range('E40').Value = '  -0.80%  '
$ws.Range('D41').Value = "'0.8498"
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('D42').Value = '1.049.57'
$ws.Range('E42').Value = '  +3.52%  '
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('D44').Value = "'102.02"
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('D45').Value = '1.887.65'
$ws.Range('E45').Value = '  +5.40%  '
$ws.Range('D46').Value = '0.0₈115'
$ws.Range('E46').Value = '  +3.78%  '
$ws.Range('D47').Value = "'59.96"
$ws.Range('E47').Value = '  +2.88%  '
$ws.Range('D48').Value = "'8.252"
$ws.Range('E48').Value = '  +2.64%  '
$ws.Range('D49').Value = "'0.4433"
$ws.Range('E49').Value = '  +2.01%  '
$ws.Range('D50').Value = "'1.003"
$ws.Range('E50').Value = '  -0.28%  '
$ws.Range('D51').Value = "'0.05324"
$ws.Range('E51').Value = '  -0.09%  '

# Reset style on cells that required a text-forcing quote prefix, so no stray
# number-format / quotePrefix styling is introduced on these cells
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').Style = "Normal"
